$wb = $excel.ActiveWorkbook

# The "BMW" sheet's contact-pair list is reused as the base data for the
# new 7.4kW CCU variants; duplicate it twice, placing each duplicate
# immediately ahead of "BMW" in tab order.
$wb.Worksheets.Item("BMW").Copy($wb.Worksheets.Item("BMW"))
$wb.Worksheets.Item("BMW (2)").Name = "7.4kW_CCU_L"

$wb.Worksheets.Item("BMW").Copy($wb.Worksheets.Item("BMW"))
$wb.Worksheets.Item("BMW (2)").Name = "7.4kW_CCU_S"

# The original "BMW" sheet becomes the 22kW CCU list.
$wb.Worksheets.Item("BMW").Name = "22kW_CCU_L"

$wb.Worksheets.Item("22kW_CCU_L").Activate()
